$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill values in the same order the cells were originally entered so the
# shared-string table comes out in the same order as the source edit
# (A column for the first two new rows, then B column for those two rows,
# then C column, then the third row filled in one go).
$ws.Range("A8").Value = "list 1"
$ws.Range("A9").Value = "list 2"

$ws.Range("B8").Value = "[a,b,c]"
$ws.Range("B9").Value = '["c","d","e"]'

$ws.Range("C8").Value = 7
$ws.Range("C9").Value = 8

$ws.Range("A10").Value = "list 3"
$ws.Range("B10").Value = "['f', 'g', 'h']"
$ws.Range("C10").Value = 9

# Match formatting of the existing data rows: column A uses the same
# style as A2/A4/A5/A6, column B uses the style used across column B.
$ws.Range("A2").Copy()
$ws.Range("A8:A10").PasteSpecial(-4122)

$ws.Range("B1").Copy()
$ws.Range("B8:B10").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# New rows keep the sheet's standard row height.
$ws.Rows.Item(8).RowHeight = 15.75
$ws.Rows.Item(9).RowHeight = 15.75
$ws.Rows.Item(10).RowHeight = 15.75

$ws.Range("B20").Select()
